$d = $word.ActiveDocument
$nl = [char]10

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Host "NOT FOUND: $old"
    }
    return $result
}

# 1) caseQuery line number 516 -> 540
Replace-Text "M2DocEvaluator.java:516)" "M2DocEvaluator.java:540)"

# 2) TemplateSwitch.java:172 -> 186
Replace-Text "TemplateSwitch.java:172)" "TemplateSwitch.java:186)"

# 3) doSwitch(M2DocEvaluator.java:945) + caseBlock(M2DocEvaluator.java:1158) -> doSwitch(1038) + caseBlock(1254)
$old3 = "org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:945)" + $nl + "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1158)"
$new3 = "org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)" + $nl + "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)"
Replace-Text $old3 $new3

# 4) TemplateSwitch.java:183 -> 199
Replace-Text "TemplateSwitch.java:183)" "TemplateSwitch.java:199)"

# 5) block: doSwitch(945) + caseTemplate(311) + caseTemplate(1) + TemplateSwitch.java:201
#    -> doSwitch(1038) + caseDocumentTemplate(275) + caseDocumentTemplate(1) + TemplateSwitch.java:279
$old5 = "org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:945)" + $nl + "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:311)" + $nl + "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseTemplate(M2DocEvaluator.java:1)" + $nl + "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:201)"
$new5 = "org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)" + $nl + "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)" + $nl + "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)" + $nl + "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)"
Replace-Text $old5 $new5

# 6) large block: doSwitch(945) + caseDocumentTemplate(266)/.../GeneratedMethodAccessor76
#    collapsing into doSwitch(1038) + generate(264)/M2DocUtils.generate(712)/prepareoutputAndGenerate(459)/generation(369)/GeneratedMethodAccessor75
$old6 = "org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:945)" + $nl + `
        "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:266)" + $nl + `
        "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)" + $nl + `
        "`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:246)" + $nl + `
        "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)" + $nl + `
        "`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)" + $nl + `
        "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:945)" + $nl + `
        "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:255)" + $nl + `
        "`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:705)" + $nl + `
        "`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:458)" + $nl + `
        "`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:365)" + $nl + `
        "`tat sun.reflect.GeneratedMethodAccessor76.invoke(Unknown Source)"

$new6 = "org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)" + $nl + `
        "`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)" + $nl + `
        "`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)" + $nl + `
        "`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)" + $nl + `
        "`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)" + $nl + `
        "`tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)"

Replace-Text $old6 $new6
